# Security-requirements slide: drop the turnaround-time / "measured on dev
# machine" bullet that used to follow "...(standalone)". We keep the
# "(standalone)" parenthetical itself and the blank lines that trailed the
# removed bullet.
#
# Text layout before the edit (each line is one run, \v == <a:br/>):
#   ...必要としない(         <- run, kept
#   スタンドアロン            <- run, kept
#   )                        <- run, REMOVED (the ")" that used to close "(standalone)")
#   \v\v                     <- 2 <a:br/>, REMOVED
#   ・予測処理の...タイムは    <- run, REMOVED
#   10                       <- run, REMOVED
#   秒以内                    <- run, REMOVED
#   \v                       <- <a:br/>, REMOVED
#     (                      <- run, REMOVED
#   ただし、...基準とする      <- run, REMOVED
#   )                        <- run, kept (becomes the new ")" closing "(standalone)")
#   \v\v                     <- 2 <a:br/>, kept

$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $txt = $shp.TextFrame.TextRange.Text
            if (($txt.IndexOf("(スタンドアロン") -ge 0) -and ($txt.IndexOf("基準とする") -ge 0)) {
                $targetSlide = $sl
                $targetShape = $shp
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$full = $tr.Text

# Anchor 1: text ending right at "...(スタンドアロン"; the char right after it
# (0-based) is the ")" we want to delete first (so the deletion range starts
# on a normal character, not on a line-break - deleting a range that *starts*
# or *ends* exactly on a <a:br/> silently leaves that break behind).
$anchor1 = "(スタンドアロン"
$idx1 = $full.IndexOf($anchor1)
$deleteStart0 = $idx1 + $anchor1.Length

# Anchor 2: text ending right at "...基準とする"; the char right after it
# (0-based) is the ")" that must survive, so the deletion stops *before* it.
$anchor2 = "基準とする"
$idx2 = $full.IndexOf($anchor2)
$deleteEnd0 = $idx2 + $anchor2.Length

$len = $deleteEnd0 - $deleteStart0

# TextRange.Characters is 1-based.
$tr.Characters($deleteStart0 + 1, $len).Delete()
